$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.873.93'
$ws.Range("E2").Value = '  +3.82%  '
$ws.Range("D3").Value = '3.098.90'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.77'
$ws.Range("E5").Value = '  +3.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.19'
$ws.Range("E6").Value = '  +7.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '3.095.66'
$ws.Range("E8").Value = '  +2.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.500'
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.42'
$ws.Range("E10").Value = '  +4.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").Value = '  +3.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.473'
$ws.Range("E12").Value = '  +6.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").Value = '  +3.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.35'
$ws.Range("E14").Value = '  +3.54%  '
$ws.Range("D15").Value = '3.601.28'
$ws.Range("E15").Value = '  +2.19%  '
$ws.Range("D16").Value = '64.901.22'
$ws.Range("E16").Value = '  +3.67%  '
$ws.Range("D17").Value = '3.100.28'
$ws.Range("E17").Value = '  +2.20%  '
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.78'
$ws.Range("E19").Value = '  +2.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.24'
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("E21").Value = '  +4.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  +2.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.57'
$ws.Range("E23").Value = '  +7.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.52'
$ws.Range("E24").Value = '  +12.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.19'
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  +3.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").Value = '  +4.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.06'
$ws.Range("E29").Value = '  +6.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.13'
$ws.Range("E31").Value = '  +1.87%  '
$ws.Range("E32").Value = '  +1.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.50'
$ws.Range("E33").Value = '  +5.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.63'
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.16'
$ws.Range("E35").Value = '  +5.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.03'
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '470.56'
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0837'
$ws.Range("E38").Value = '  +4.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0408'
$ws.Range("E39").Value = '  +5.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.94'
$ws.Range("E40").Value = '  +20.31%  '
$ws.Range("D41").Value = '2.982.45'
$ws.Range("E41").Value = '  -5.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.26'
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.30'
$ws.Range("E44").Value = '  +6.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.260'
$ws.Range("E45").Value = '  +5.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.15'
$ws.Range("E46").Value = '  +8.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.112'
$ws.Range("E48").Value = '  +3.49%  '
$ws.Range("D49").Value = '0.0₃0530'
$ws.Range("E49").Value = '  +5.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.74'
$ws.Range("E50").Value = '  +2.72%  '
$ws.Range("E51").Value = '  +3.28%  '
